$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 = "Save", formatted like the other header cells (e.g. G1 = bold/border style)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# New data cell H2 = 0 (numeric)
$ws.Range("H2").Value = 0
